$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Row 2 (condition = light-dark): mean/sd/n/se/cv
$table.Cell(2, 1).Range.Text = "suspension"
$table.Cell(2, 2).Range.Text = "61.82"
$table.Cell(2, 3).Range.Text = "28.53"
$table.Cell(2, 4).Range.Text = "136"
$table.Cell(2, 5).Range.Text = "2.45"
$table.Cell(2, 6).Range.Text = "3.96"

# Row 3 (condition = light): mean/sd/n/se/cv
$table.Cell(3, 1).Range.Text = "plug"
$table.Cell(3, 2).Range.Text = "41.98"
$table.Cell(3, 3).Range.Text = "31.50"
$table.Cell(3, 4).Range.Text = "128"
$table.Cell(3, 5).Range.Text = "2.78"
$table.Cell(3, 6).Range.Text = "6.63"
